# Update country case/death statistics and reorder two pairs of
# countries (Timor Oriental <-> Santa Lucia, Islas Malvinas <-> Montserrat)
# whose ranking changed places, plus refresh the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Per-country case/death figures refreshed for this snapshot ---------
# Most countries simply have "Nuevos casos" (C) and "Muertes hoy" (G) reset
# to 0 (today's counters cleared); a handful of countries also received
# updated cumulative totals (B/D/E/H) for this run.
$rowUpdates = @(
    @{Row=4; B=6925941; C=0; D=4191894; E=2530876; G=0; H=203171}
    @{Row=5; C=0; G=0}
    @{Row=6; C=0; G=0}
    @{Row=7; C=0; G=0}
    @{Row=8; C=0; G=0}
    @{Row=9; C=0; G=0}
    @{Row=10; B=688954; C=4841; D=492192; E=123959; G=624; H=72803}
    @{Row=11; C=0; G=0}
    @{Row=12; C=0; G=0}
    @{Row=13; C=0; G=0}
    @{Row=14; C=0; G=0}
    @{Row=15; C=0; G=0}
    @{Row=16; C=0; G=0}
    @{Row=17; C=0; G=0}
    @{Row=18; C=0; G=0}
    @{Row=19; C=0; G=0}
    @{Row=20; C=0; G=0}
    @{Row=21; C=0; G=0}
    @{Row=22; C=0; G=0}
    @{Row=23; C=0; G=0}
    @{Row=24; C=0; G=0}
    @{Row=25; C=0; D=243000; E=18780; G=0}
    @{Row=26; C=0; G=0}
    @{Row=27; C=0; G=0}
    @{Row=28; C=0; G=0}
    @{Row=29; C=0; G=0}
    @{Row=30; B=130051; C=632; D=88457; E=34044; G=39; H=7550}
    @{Row=31; C=0; G=0}
    @{Row=32; C=0; G=0}
    @{Row=33; C=0; G=0}
    @{Row=34; C=0}
    @{Row=35; C=0; G=0}
    @{Row=36; C=0; G=0}
    @{Row=37; C=0; G=0}
    @{Row=38; C=0; G=0}
    @{Row=39; C=0; G=0}
    @{Row=40; C=0; G=0}
    @{Row=42; C=0; G=0}
    @{Row=43; G=0}
    @{Row=44; B=85269; C=14; D=80464; E=171}
    @{Row=45; C=0; G=0}
    @{Row=46; C=0; G=0}
    @{Row=47; C=0; G=0}
    @{Row=48; C=0; G=0}
    @{Row=49; C=0; G=0}
    @{Row=50; C=0; G=0}
    @{Row=51; C=0; G=0}
    @{Row=52; C=0; G=0}
    @{Row=54; C=0; G=0}
    @{Row=55; C=0; G=0}
    @{Row=56; C=0; G=0}
    @{Row=57; C=0}
    @{Row=58; C=0; G=0}
    @{Row=59; C=0; G=0}
    @{Row=60; C=0; G=0}
    @{Row=61; C=0; G=0}
    @{Row=62; C=0; G=0}
    @{Row=63; C=0; G=0}
    @{Row=64; C=0; G=0}
    @{Row=65; C=0; G=0}
    @{Row=66; C=0}
    @{Row=67; C=0; G=0}
    @{Row=68; C=0; G=0}
    @{Row=69; C=0; G=0}
    @{Row=70; C=0; G=0}
    @{Row=71; C=0; G=0}
    @{Row=72; C=0; G=0}
    @{Row=73; C=0; G=0}
    @{Row=74; C=0; G=0}
    @{Row=75; C=0; G=0}
    @{Row=76; C=0; G=0}
    @{Row=77; B=26885; C=24; D=23861; E=2180; G=7; H=844}
    @{Row=78; C=0; G=0}
    @{Row=79; C=0; G=0}
    @{Row=80; B=22893; C=110; D=19970; E=2545; G=1; H=378}
    @{Row=81; C=0}
    @{Row=82; C=0; G=0}
    @{Row=83; C=0}
    @{Row=84; C=0; G=0}
    @{Row=85; C=0; G=0}
    @{Row=86; C=0; G=0}
    @{Row=87; C=0; G=0}
    @{Row=88; C=0; G=0}
    @{Row=89; C=0; G=0}
    @{Row=90; C=0; G=0}
    @{Row=91; C=0; G=0}
    @{Row=93; C=0; G=0}
    @{Row=94; C=0; G=0}
    @{Row=95; C=0; G=0}
    @{Row=96; C=0}
    @{Row=97; C=0}
    @{Row=98; C=0; G=0}
    @{Row=99; C=0}
    @{Row=100; C=0}
    @{Row=101; C=0}
    @{Row=102; C=0}
    @{Row=103; C=0}
    @{Row=104; C=0; G=0}
    @{Row=106; C=0}
    @{Row=107; C=0; G=0}
    @{Row=108; C=0}
    @{Row=109; C=0}
    @{Row=110; C=0; G=0}
    @{Row=111; C=0}
    @{Row=112; C=0}
    @{Row=113; C=0; G=0}
    @{Row=114; C=0}
    @{Row=115; C=0}
    @{Row=116; C=0; G=0}
    @{Row=117; C=0; G=0}
    @{Row=119; C=0; G=0}
    @{Row=120; C=0; G=0}
    @{Row=122; C=0}
    @{Row=123; C=0; G=0}
    @{Row=124; C=0; G=0}
    @{Row=125; C=0; G=0}
    @{Row=126; C=0; G=0}
    @{Row=127; C=0; G=0}
    @{Row=128; C=0; G=0}
    @{Row=129; C=0; G=0}
    @{Row=130; C=0; G=0}
    @{Row=131; C=0}
    @{Row=132; C=0}
    @{Row=134; C=0}
    @{Row=135; C=0; G=0}
    @{Row=136; C=0}
    @{Row=139; C=0}
    @{Row=140; C=0}
    @{Row=142; C=0}
    @{Row=143; C=0}
    @{Row=144; C=0}
    @{Row=145; C=0; G=0}
    @{Row=146; C=0}
    @{Row=148; C=0}
    @{Row=150; C=0}
    @{Row=151; C=0}
    @{Row=152; C=0; G=0}
    @{Row=153; C=0}
    @{Row=154; C=0}
    @{Row=155; B=1811; C=2; D=1719; E=67}
    @{Row=156; C=0}
    @{Row=157; C=0}
    @{Row=158; C=0}
    @{Row=159; C=0}
    @{Row=160; C=0}
    @{Row=161; C=0}
    @{Row=162; C=0}
    @{Row=163; C=0}
    @{Row=165; C=0}
    @{Row=168; C=0}
    @{Row=169; C=0}
    @{Row=172; C=0}
    @{Row=173; C=0}
    @{Row=179; C=0}
    @{Row=182; C=0}
    @{Row=187; C=0}
    @{Row=188; C=0}
    @{Row=190; C=0}
    @{Row=201; C=0}
    @{Row=214; D=13; H=0}
    @{Row=215; D=12; H=1}
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey('D')) { $ws.Cells.Item($r, 4).Value = $u.D }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
    if ($u.ContainsKey('F')) { $ws.Cells.Item($r, 6).Value = $u.F }
    if ($u.ContainsKey('G')) { $ws.Cells.Item($r, 7).Value = $u.G }
    if ($u.ContainsKey('H')) { $ws.Cells.Item($r, 8).Value = $u.H }
}

# --- Two countries swapped ranking position in the source table ---------
# "Timor Oriental" now ranks above "Santa Lucia" (rows 204/205), and
# "Islas Malvinas" now ranks above "Montserrat" (rows 214/215). Swap the
# country names between the two fixed rows.
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# --- Refresh the "last updated" timestamp --------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 04:11"
